# Fruta / hortaliza, semanal
# Insert a new daily price record as row 343 (pushing the existing rows
# 343..452 down to 344..453), matching the new market-day entry added to
# the "Fruta, Feria Lagunitas de Puerto Montt - Mango" consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 343 down by one row (Excel's normal
# "insert row" behaviour - carries down formats, dimension grows to T453).
$ws.Rows("343:343").Insert()

# Populate the newly-blank row 343 with the new record.
$ws.Range("A343").Value = 4
$ws.Range("B343").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C343").Value = "Los Lagos"
$ws.Range("D343").Value = 45215
$ws.Range("E343").Value = 10
$ws.Range("F343").Value = "Fruta"
$ws.Range("G343").Value = 100108
$ws.Range("H343").Value = "Tropicales y subtropicales"
$ws.Range("I343").Value = 100108002
$ws.Range("J343").Value = "Mango"
$ws.Range("K343").Value = "Sin especificar"
$ws.Range("L343").Value = "Primera"
$ws.Range("M343").Value = 100
$ws.Range("N343").Value = 13000
$ws.Range("O343").Value = 13000
$ws.Range("P343").Value = 13000
$ws.Range("Q343").Value = "`$/bandeja 4 kilos"
$ws.Range("R343").Value = "Brasil"
$ws.Range("S343").Value = 3250
$ws.Range("T343").Value = 4
